# LOB1215.docx edit script
# Applies the changes described by the commit diff using Word COM-interop
# Find/Replace calls, scoped to specific paragraphs where the "before"
# text is not unique within the document (so a document-wide Replace-All
# would otherwise touch the wrong paragraph too).

$d = $word.ActiveDocument
$vt = [char]11   # vertical-tab char == <w:br/> manual line break in Range.Text

function Replace-InRange($range, [string]$findText, [string]$replaceText) {
    $ok = $range.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

# 1) Update activation date (unique string, safe to do document-wide)
Replace-InRange $d.Content "Ativação: 01/01/2015" "Ativação: 01/01/2025"

# 2) Remove the second "responsible professor" (Wendell), merging the
#    break + run into the first run so only Ana Karine remains.
$findWendell = "7043088 - Ana Karine Furtado de Carvalho" + $vt + "8426375 - Wendell de Queiróz Lamas"
Replace-InRange $d.Content $findWendell "7043088 - Ana Karine Furtado de Carvalho"

# 3) "Programa resumido" section (Portuguese) - paragraph 11
$findResumoPt = "Recursos energéticos e matriz energética do Brasil. Atividades antrópicas," + $vt + "demanda de energia e desenvolvimento socioeconômico. Disponibilidade de fontes e avaliação do potencial de geração de energia. Energia elétrica: fundamentos sobre geração, transmissão e distribuição. Usinas hidroelétricas, termoelétricas e nucleares. Energia solar. Energia eólica. Energia fóssil. Energia da biomassa. Impactos ambientais decorrentes da geração, transmissão, disponibilidade e oferta de energia no desenvolvimento regional."
$newResumoPt = "Recursos energéticos e matriz energética do Brasil. Atividades antrópicas,demanda de energia e desenvolvimento socioeconômico. Disponibilidade de fontes e avaliação do potencial de geração de energia. Energia elétrica: fundamentos sobre geração, transmissão e distribuição. Usinas hidroelétricas, termoelétricas e nucleares. Energia fóssil. Impactos ambientais decorrentes da geração, transmissão, disponibilidade e oferta de energia no desenvolvimento regional."
Replace-InRange $d.Paragraphs(11).Range $findResumoPt $newResumoPt

# 4) "Programa resumido" section (English / italic) - paragraph 12
$findResumoEn = "Energy sources and the Brazilian energy matrix. Anthropogenic activities, energy demand, and socio-economic development. Availability of sources and evaluation of energy generation potential. Electrical power fundamentals. Power plants. Solar energy. Wind energy. Fossil energy. Biomass energy. Environmental impacts of energy generation, transmission, availability, and supply in regional development."
$newResumoEn = "Energy sources and the Brazilian energy matrix. Anthropogenic activities, energy demand, and socio-economic development. Availability of sources and evaluation of energy generation potential. Electrical power fundamentals. Power plants. Fossil energy. Environmental impacts of energy generation, transmission, availability, and supply in regional development."
Replace-InRange $d.Paragraphs(12).Range $findResumoEn $newResumoEn

# 5) "Programa" section (Portuguese) - paragraph 14
$findProgPt = "Recursos energéticos e matriz energética do Brasil. Atividades antrópicas," + $vt + "demanda de energia e desenvolvimento socioeconômico. Disponibilidade de fontes e avaliação do potencial de geração de energia. Energia elétrica: fundamentos sobre geração, transmissão e distribuição. Usinas hidroelétricas, termoelétricas e nucleares. Energia solar. Energia eólica. Energia fóssil. Energia da biomassa. Impactos ambientais decorrentes da geração, transmissão, disponibilidade e oferta de energia no desenvolvimento regional."
$newProgPt = "Recursos energéticos e matriz energética do Brasil. Atividades antrópicas,demanda de energia e desenvolvimento socioeconômico. Disponibilidade de fontes e avaliação do potencial de geração de energia. Energia elétrica: fundamentos sobre geração, transmissão e distribuição. Usinas hidroelétricas, termoelétricas e nucleares. Energia fóssil. Impactos ambientais decorrentes da geração, transmissão, disponibilidade e oferta de energia no desenvolvimento regiona. A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina."
Replace-InRange $d.Paragraphs(14).Range $findProgPt $newProgPt

# 6) "Programa" section (English / italic) - paragraph 15
$findProgEn = "Energy sources and the Brazilian energy matrix. Anthropogenic activities, energy demand, and socio-economic development. Availability of sources and evaluation of energy generation potential. Electrical power fundamentals: generation, transmission, and distribution. Power plants: hydraulic, thermal, and nuclear. Solar energy. Wind energy. Fossil energy. Biomass energy. Environmental impacts of energy generation, transmission, availability, and supply in regional development."
$newProgEn = "Energy sources and the Brazilian energy matrix. Anthropogenic activities, energy demand, and socio-economic development. Availability of sources and evaluation of energy generation potential. Electrical power fundamentals: generation, transmission, and distribution. Power plants: hydraulic, thermal, and nuclear. Fossil energy. Environmental impacts of energy generation, transmission, availability, and supply in regional development. The discipline may have didactic trips to complement the content of the discipline."
Replace-InRange $d.Paragraphs(15).Range $findProgEn $newProgEn

# 7) Avaliação - Método
Replace-InRange $d.Content "Os alunos efetuarão monografias em grupos a serem selecionados em classe." "O método de avaliação será composto por avaliação teórica, apresentação escrita e oral."

# 8) Avaliação - Critério
Replace-InRange $d.Content "Dois seminários, pesos 1 e 2." "Para o cálculo da nota final (NF) será adotada a média ponderada de provas e atividades."

# 9) Avaliação - Norma de recuperação
Replace-InRange $d.Content "Para os alunos reprovados por nota, mas beneficiados pelo sistema de recuperação, esta será realizada através da aplicação de uma única prova teórica, abrangendo todo o programa do semestre letivo." "Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 => 5,0 Aprovado"

Write-Output "done"
